$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new source entry (fix for emoji handling) as row 9
$ws.Range("A9").Value = "emojis_unicode"
$ws.Range("B9").Value = "dictionary"
$ws.Range("C9").Value = "https://github.com/today-is-a-good-day/emojis/blob/master/emDict.csv"

# Update selection/view to the newly added cell
$ws.Range("C9").Select()
